$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values that look numeric stay as text (matches original inlineStr formatting)
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '37.993.83'
$ws.Range('E2').Value = '  +2.91%  '
$ws.Range('D3').Value = '2.054.47'
$ws.Range('E3').Value = '  +2.35%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = '230.20'
$ws.Range('E5').Value = '  +2.14%  '
$ws.Range('D6').Value = '0.616'
$ws.Range('E6').Value = '  +2.08%  '
$ws.Range('E7').Value = '  +8.06%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +3.44%  '
$ws.Range('E10').Value = '  +4.93%  '
$ws.Range('D11').Value = '0.103'
$ws.Range('E11').Value = '  +1.92%  '
$ws.Range('D12').Value = '2.358.44'
$ws.Range('E12').Value = '  +2.15%  '
$ws.Range('D13').Value = '14.64'
$ws.Range('E13').Value = '  +4.83%  '
$ws.Range('E14').Value = '  +5.65%  '
$ws.Range('D15').Value = '0.753'
$ws.Range('E15').Value = '  +2.71%  '
$ws.Range('E16').Value = '  +1.76%  '
$ws.Range('D17').Value = '2.056.08'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '37.896.06'
$ws.Range('E18').Value = '  +2.85%  '
$ws.Range('D19').Value = '6.37'
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('E20').Value = '  +2.30%  '
$ws.Range('D21').Value = '0.0₃0837'
$ws.Range('E21').Value = '  +3.51%  '
$ws.Range('D22').Value = '224.63'
$ws.Range('E22').Value = '  +1.47%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').Value = '2.44'
$ws.Range('E24').Value = '  +0.19%  '
$ws.Range('D25').Value = '2.25'
$ws.Range('E25').Value = '  +4.48%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '166.73'
$ws.Range('E26').Value = '  +1.17%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '9.29'
$ws.Range('E27').Value = '  +2.62%  '
$ws.Range('E28').Value = '  +5.14%  '
$ws.Range('E29').Value = '  +2.95%  '
$ws.Range('E30').Value = '  +3.07%  '
$ws.Range('E31').Value = '  +3.18%  '
$ws.Range('D32').Value = '4.54'
$ws.Range('E32').Value = '  +2.13%  '
$ws.Range('D33').Value = '4.60'
$ws.Range('E33').Value = '  +4.15%  '
$ws.Range('E34').Value = '  +10.78%  '
$ws.Range('E35').Value = '  +1.82%  '
$ws.Range('E36').Value = '  +1.06%  '
$ws.Range('D37').Value = '6.06'
$ws.Range('E37').Value = '  +13.20%  '
$ws.Range('D38').Value = '3.28'
$ws.Range('E38').Value = '  +6.67%  '
$ws.Range('E39').Value = '  -0.32%  '
$ws.Range('D40').Value = '1.507.69'
$ws.Range('E40').Value = '  +3.81%  '
$ws.Range('E41').Value = '  +3.31%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').Value = '2.89'
$ws.Range('E42').Value = '  +4.27%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '97.02'
$ws.Range('E43').Value = '  +2.66%  '
$ws.Range('D44').Value = '16.54'
$ws.Range('E44').Value = '  +3.50%  '
$ws.Range('E45').Value = '  +1.37%  '
$ws.Range('E46').Value = '  +0.74%  '
$ws.Range('D47').Value = '4.11'
$ws.Range('E47').Value = '  +16.13%  '
$ws.Range('E48').Value = '  +2.23%  '
$ws.Range('E49').Value = '  +2.24%  '
$ws.Range('D50').Value = '7.10'
$ws.Range('E50').Value = '  -0.36%  '
$ws.Range('D51').Value = '2.245.75'
$ws.Range('E51').Value = '  +2.16%  '
